# Fix bug: rows for "Expand Security Council..." and "International levy on
# shipping carbon emissions..." were swapped, and the "All" / "Russia"
# aggregate columns needed recalculating.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the row labels in column A for rows 5 and 6 ---
$ws.Range("A5").Value = "International levy on shipping carbon emissions,`nreturned to countries based on population"
$ws.Range("A6").Value = "Expand Security Council to new permanent members (e.g.`nIndia, Brazil, African Union), restrict veto use"

# --- Row 2 ---
$ws.Range("B2").Value = 0.807885696236268
$ws.Range("L2").Value = 0.790720704052685

# --- Row 3 ---
$ws.Range("B3").Value = 0.794189426632509
$ws.Range("L3").Value = 0.83450478996974

# --- Row 4 ---
$ws.Range("B4").Value = 0.74891324454649
$ws.Range("L4").Value = 0.873558433686924

# --- Row 5 (now "International levy on shipping...") ---
$ws.Range("B5").Value = 0.700357100296232
$ws.Range("C5").Value = 0.732870289440397
$ws.Range("D5").Value = 0.776964615344406
$ws.Range("E5").Value = 0.696966943798898
$ws.Range("F5").Value = 0.779538293526919
$ws.Range("G5").Value = 0.612876561840658
$ws.Range("H5").Value = 0.738317533612799
$ws.Range("I5").Value = 0.747513876364388
$ws.Range("J5").Value = 0.718133301168214
$ws.Range("K5").Value = 0.58759070342811
$ws.Range("L5").Value = 0.731306566882626
$ws.Range("M5").Value = 0.814900578705803
$ws.Range("N5").Value = 0.674291226582879

# --- Row 6 (now "Expand Security Council...") ---
$ws.Range("B6").Value = 0.698034537168434
$ws.Range("C6").Value = 0.761932435910776
$ws.Range("D6").Value = 0.720303245818442
$ws.Range("E6").Value = 0.756868069241599
$ws.Range("F6").Value = 0.802386607493625
$ws.Range("G6").Value = 0.7250772374002
$ws.Range("H6").Value = 0.763310673110362
$ws.Range("I6").Value = 0.78061912302082
$ws.Range("J6").Value = 0.721726707430624
$ws.Range("K6").Value = 0.67685065346556
$ws.Range("L6").Value = 0.552689160628133
$ws.Range("M6").Value = 0.836791104476278
$ws.Range("N6").Value = 0.671916850954642

# --- Row 7 ---
$ws.Range("B7").Value = 0.696942403411651
$ws.Range("L7").Value = 0.813558582458709

# --- Row 8 ---
$ws.Range("B8").Value = 0.694726686848076
$ws.Range("C8").Value = 0.698508510826338
$ws.Range("L8").Value = 0.730009595228465

# --- Row 9 ---
$ws.Range("B9").Value = 0.682188628877153
$ws.Range("L9").Value = 0.868962834674846

# --- Row 10 ---
$ws.Range("B10").Value = 0.68030147177849
$ws.Range("L10").Value = 0.493073564956229

# --- Row 11 ---
$ws.Range("B11").Value = 0.527237118011288
$ws.Range("L11").Value = 0.519297451909645
